# 1st changes of mifos to finflux
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted right
# before the existing "Late" column (column N), pushing the "Late" and
# "Outstanding" columns one place to the right (N->O, O->P, P->Q).
#
# The active/selected sheet also switches from "Transactions" to
# "Repayment Schedule", with a new selected cell of T4 on that sheet.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q).
$wsSchedule.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with T4 selected,
# which also clears the tabSelected/active state on "Transactions".
$wsSchedule.Activate()
$wsSchedule.Range("T4").Select()
